$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update servo trim values in the new column K ("trim") for rows 3-14
$ws.Range("K3").Value = -6
$ws.Range("K4").Value = -9
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = -7
$ws.Range("K7").Value = 6
$ws.Range("K8").Value = 0
$ws.Range("K9").Value = 54
$ws.Range("K10").Value = -44
$ws.Range("K11").Value = -3
$ws.Range("K12").Value = -5
$ws.Range("K13").Value = 32
$ws.Range("K14").Value = 9

# Correct the pin value for servo15 (row 15), then add its trim value
$ws.Range("J15").Value = 26
$ws.Range("K15").Value = 3

# Continue filling in trim values for the remaining rows
$ws.Range("K16").Value = -6
$ws.Range("K17").Value = -2
$ws.Range("K18").Value = -4
$ws.Range("K19").Value = -2
$ws.Range("K20").Value = 2
$ws.Range("K21").Value = 2
$ws.Range("K22").Value = 0
$ws.Range("K23").Value = -5
$ws.Range("K24").Value = 2

# Update the saved view state (scroll position & active selection)
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("K25").Select()
